# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical event listings, and the same set of rows
# received updated counts in this edit.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 60
    5  = 76
    8  = 1496
    10 = 371
    12 = 134
    14 = 61
    15 = 102
    17 = 290
    19 = 1700
    21 = 106
    22 = 171
    23 = 646
    25 = 330
    26 = 4089
    28 = 476
    29 = 254
    30 = 1057
    31 = 129
    33 = 418
    35 = 186
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
